# Applies the "Interpretación mejor modelo MR" update to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Ridge
$ws.Range("B2").Value = 0.9859789073419475
$ws.Range("C2").Value = 0.9235314557650617
$ws.Range("D2").Value = 0.0624474515768858
$ws.Range("E2").Value = 3.883917171951741
$ws.Range("F2").Value = 9.070274720007825
$ws.Range("H2").Value = 3.87

# Row 3 - Lasso
$ws.Range("B3").Value = 0.986294233370364
$ws.Range("C3").Value = 0.9244396652579883
$ws.Range("D3").Value = 0.06185456811237577
$ws.Range("E3").Value = 3.839995330965513
$ws.Range("F3").Value = 9.016250443449955
$ws.Range("G3").Value = "{'selection': 'cyclic', 'alpha': 0.01}"
$ws.Range("H3").Value = 0.74

# Row 4 - ElasticNet
$ws.Range("H4").Value = 1.99

# Row 5 - SVR
$ws.Range("H5").Value = 4.98

# Row 6 - KNN Regressor
$ws.Range("H6").Value = 5.72

# Row 7 - Decision Tree
$ws.Range("H7").Value = 2.05

# Row 8 - PLSRegression
$ws.Range("H8").Value = 0.5
